# Refresh the cryptocurrency Price (D) and Volume(1h) (E) columns for rows 2-51.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.263.00"
$ws.Range("E2").Value = "  -0.55%  "
$ws.Range("D3").Value = "1.805.86"
$ws.Range("E3").Value = "  -0.66%  "
$ws.Range("E4").Value = "  +0.07%  "
$r = $ws.Range("D5")
$r.NumberFormat = "@"
$r.Value = "314.62"
$r.ClearFormats()
$ws.Range("E5").Value = "  -0.18%  "
$r = $ws.Range("D6")
$r.NumberFormat = "@"
$r.Value = "1.002"
$r.ClearFormats()
$ws.Range("E6").Value = "  +0.05%  "
$r = $ws.Range("D7")
$r.NumberFormat = "@"
$r.Value = "0.5267"
$r.ClearFormats()
$ws.Range("E7").Value = "  +2.68%  "
$r = $ws.Range("D8")
$r.NumberFormat = "@"
$r.Value = "0.3827"
$r.ClearFormats()
$ws.Range("E8").Value = "  -3.06%  "
$r = $ws.Range("D9")
$r.NumberFormat = "@"
$r.Value = "0.08033"
$r.ClearFormats()
$ws.Range("E9").Value = "  +0.40%  "
$ws.Range("E10").Value = "  -0.65%  "
$r = $ws.Range("D11")
$r.NumberFormat = "@"
$r.Value = "1.103"
$r.ClearFormats()
$ws.Range("E11").Value = "  -0.51%  "
$r = $ws.Range("D12")
$r.NumberFormat = "@"
$r.Value = "6.336"
$r.ClearFormats()
$ws.Range("E12").Value = "  +1.28%  "
$ws.Range("E13").Value = "  +0.10%  "
$r = $ws.Range("D14")
$r.NumberFormat = "@"
$r.Value = "20.65"
$r.ClearFormats()
$ws.Range("E14").Value = "  -1.57%  "
$r = $ws.Range("D15")
$r.NumberFormat = "@"
$r.Value = "7.339"
$r.ClearFormats()
$ws.Range("E15").Value = "  -2.06%  "
$ws.Range("D16").Value = "1.806.40"
$ws.Range("E16").Value = "  -1.45%  "
$r = $ws.Range("D17")
$r.NumberFormat = "@"
$r.Value = "92.30"
$r.ClearFormats()
$ws.Range("E17").Value = "  -0.39%  "
$r = $ws.Range("D18")
$r.NumberFormat = "@"
$r.Value = "0.00001098"
$r.ClearFormats()
$ws.Range("E18").Value = "  -3.11%  "
$r = $ws.Range("D19")
$r.NumberFormat = "@"
$r.Value = "0.06610"
$r.ClearFormats()
$ws.Range("E19").Value = "  -0.37%  "
$ws.Range("E20").Value = "  +0.07%  "
$ws.Range("E21").Value = "  -1.45%  "
$r = $ws.Range("D22")
$r.NumberFormat = "@"
$r.Value = "5.973"
$r.ClearFormats()
$ws.Range("E22").Value = "  -1.95%  "
$ws.Range("D23").Value = "28.319.25"
$ws.Range("E23").Value = "  -0.45%  "
$r = $ws.Range("D24")
$r.NumberFormat = "@"
$r.Value = "11.17"
$r.ClearFormats()
$ws.Range("E24").Value = "  -0.75%  "
$r = $ws.Range("D25")
$r.NumberFormat = "@"
$r.Value = "2.256"
$r.ClearFormats()
$ws.Range("E25").Value = "  -0.60%  "
$r = $ws.Range("D26")
$r.NumberFormat = "@"
$r.Value = "160.89"
$r.ClearFormats()
$ws.Range("E26").Value = "  +3.45%  "
$r = $ws.Range("D27")
$r.NumberFormat = "@"
$r.Value = "20.49"
$r.ClearFormats()
$ws.Range("E27").Value = "  -2.92%  "
$ws.Range("D28").Value = "2.010.34"
$ws.Range("E28").Value = "  -1.28%  "
$r = $ws.Range("D29")
$r.NumberFormat = "@"
$r.Value = "2.368"
$r.ClearFormats()
$ws.Range("E29").Value = "  -1.50%  "
$r = $ws.Range("D30")
$r.NumberFormat = "@"
$r.Value = "123.40"
$r.ClearFormats()
$ws.Range("E30").Value = "  -1.86%  "
$r = $ws.Range("D31")
$r.NumberFormat = "@"
$r.Value = "0.1085"
$r.ClearFormats()
$ws.Range("E31").Value = "  -1.51%  "
$r = $ws.Range("D32")
$r.NumberFormat = "@"
$r.Value = "1.057"
$r.ClearFormats()
$ws.Range("E32").Value = "  -4.04%  "
$ws.Range("E33").Value = "  +0.97%  "
$r = $ws.Range("D34")
$r.NumberFormat = "@"
$r.Value = "5.568"
$r.ClearFormats()
$ws.Range("E34").Value = "  -2.52%  "
$r = $ws.Range("D35")
$r.NumberFormat = "@"
$r.Value = "0.07294"
$r.ClearFormats()
$ws.Range("E35").Value = "  +3.90%  "
$r = $ws.Range("D36")
$r.NumberFormat = "@"
$r.Value = "12.49"
$r.ClearFormats()
$ws.Range("E36").Value = "  +10.71%  "
$r = $ws.Range("D37")
$r.NumberFormat = "@"
$r.Value = "0.02326"
$r.ClearFormats()
$ws.Range("E37").Value = "  +0.04%  "
$ws.Range("E38").Value = "  -2.70%  "
$r = $ws.Range("D39")
$r.NumberFormat = "@"
$r.Value = "5.150"
$r.ClearFormats()
$ws.Range("E39").Value = "  -0.85%  "
$r = $ws.Range("D40")
$r.NumberFormat = "@"
$r.Value = "8.661"
$r.ClearFormats()
$ws.Range("E40").Value = "  -1.90%  "
$r = $ws.Range("D41")
$r.NumberFormat = "@"
$r.Value = "0.6222"
$r.ClearFormats()
$ws.Range("E41").Value = "  -0.50%  "
$ws.Range("E42").Value = "  -0.67%  "
$r = $ws.Range("D43")
$r.NumberFormat = "@"
$r.Value = "1.370"
$r.ClearFormats()
$ws.Range("E43").Value = "  -2.19%  "
$r = $ws.Range("D44")
$r.NumberFormat = "@"
$r.Value = "0.6046"
$r.ClearFormats()
$ws.Range("E44").Value = "  +2.37%  "
$r = $ws.Range("D45")
$r.NumberFormat = "@"
$r.Value = "13.18"
$r.ClearFormats()
$ws.Range("E45").Value = "  -1.84%  "
$r = $ws.Range("D47")
$r.NumberFormat = "@"
$r.Value = "127.30"
$r.ClearFormats()
$ws.Range("E47").Value = "  +1.99%  "
$r = $ws.Range("D48")
$r.NumberFormat = "@"
$r.Value = "1.220"
$r.ClearFormats()
$ws.Range("E48").Value = "  +2.70%  "
$r = $ws.Range("D49")
$r.NumberFormat = "@"
$r.Value = "1.933"
$r.ClearFormats()
$ws.Range("E49").Value = "  -2.01%  "
$r = $ws.Range("D50")
$r.NumberFormat = "@"
$r.Value = "0.06817"
$r.ClearFormats()
$ws.Range("E50").Value = "  -1.02%  "
$r = $ws.Range("D51")
$r.NumberFormat = "@"
$r.Value = "73.25"
$r.ClearFormats()
$ws.Range("E51").Value = "  -1.61%  "
